# "check all chapt and add thanks, abstractE, personal publication"
#
# Record that the first draft of the thesis is done: row 7 gets a date
# (2016-04-02) in column A and the note "论文初稿完成" in column B, then
# move the active selection down to A8 (the next empty row), matching
# where the author continued working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A7").Value = 42462
$ws.Range("B7").Value = "论文初稿完成"

$ws.Range("A8").Select()
